$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: M1 -> Ccl4 -> Ccr1 -> M1 -------------------------------------
$ws.Range("B2").Value = "Ccl4"
$ws.Range("C2").Value = "Ccr1"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 64.472663
$ws.Range("H2").Value = 193.417989
$ws.Range("I2").Value = 0.4634539850319053
$ws.Range("J2").Value = 0.4634539850319053
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 58.378819
$ws.Range("N2").Value = 175.136457
$ws.Range("O2").Value = 0.6920327730022572
$ws.Range("P2").Value = 0.6920327730022573
$ws.Range("Q2").Value = 3763.837923724997
$ws.Range("R2").Value = 33874.54131352498
$ws.Range("S2").Value = 0.320725346420576
$ws.Range("T2").Value = 0.320725346420576

# --- Row 3: M1 -> Ccl4 -> Ccr1 -> M2 -------------------------------------
$ws.Range("B3").Value = "Ccl4"
$ws.Range("C3").Value = "Ccr1"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 64.472663
$ws.Range("H3").Value = 193.417989
$ws.Range("I3").Value = 0.4634539850319053
$ws.Range("J3").Value = 0.4634539850319053
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 25.979641
$ws.Range("N3").Value = 77.938923
$ws.Range("O3").Value = 0.3079672269977427
$ws.Range("P3").Value = 0.3079672269977428
$ws.Range("Q3").Value = 1674.976639053983
$ws.Range("R3").Value = 15074.78975148585
$ws.Range("S3").Value = 0.1427286386113292
$ws.Range("T3").Value = 0.1427286386113293

# --- Row 4: M2 -> Ccl4 -> Ccr1 -> M1 -------------------------------------
$ws.Range("B4").Value = "Ccl4"
$ws.Range("C4").Value = "Ccr1"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 74.50600166666666
$ws.Range("H4").Value = 223.518005
$ws.Range("I4").Value = 0.5355774335117884
$ws.Range("J4").Value = 0.5355774335117883
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 58.378819
$ws.Range("N4").Value = 175.136457
$ws.Range("O4").Value = 0.6920327730022572
$ws.Range("P4").Value = 0.6920327730022573
$ws.Range("Q4").Value = 4349.572385712031
$ws.Range("R4").Value = 39146.15147140829
$ws.Range("S4").Value = 0.370637136470595
$ws.Range("T4").Value = 0.370637136470595

# --- Row 5: M2 -> Ccl4 -> Ccr1 -> M2 -------------------------------------
$ws.Range("B5").Value = "Ccl4"
$ws.Range("C5").Value = "Ccr1"
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 74.50600166666666
$ws.Range("H5").Value = 223.518005
$ws.Range("I5").Value = 0.5355774335117884
$ws.Range("J5").Value = 0.5355774335117883
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.979641
$ws.Range("N5").Value = 77.938923
$ws.Range("O5").Value = 0.3079672269977427
$ws.Range("P5").Value = 0.3079672269977428
$ws.Range("Q5").Value = 1935.639175645402
$ws.Range("R5").Value = 17420.75258080861
$ws.Range("S5").Value = 0.1649402970411934
$ws.Range("T5").Value = 0.1649402970411934

# --- Row 6 (new): Neutro -> Ccl4 -> Ccr1 -> M1 ---------------------------
$ws.Range("A6").Value = "Neutro"
$ws.Range("B6").Value = "Ccl4"
$ws.Range("C6").Value = "Ccr1"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1347426666666666
$ws.Range("H6").Value = 0.404228
$ws.Range("I6").Value = 0.000968581456306409
$ws.Range("J6").Value = 0.0009685814563064089
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 58.378819
$ws.Range("N6").Value = 175.136457
$ws.Range("O6").Value = 0.6920327730022572
$ws.Range("P6").Value = 0.6920327730022573
$ws.Range("Q6").Value = 7.866117748910666
$ws.Range("R6").Value = 70.795059740196
$ws.Range("S6").Value = 0.0006702901110862889
$ws.Range("T6").Value = 0.0006702901110862889

# --- Row 7 (new): Neutro -> Ccl4 -> Ccr1 -> M2 ---------------------------
$ws.Range("A7").Value = "Neutro"
$ws.Range("B7").Value = "Ccl4"
$ws.Range("C7").Value = "Ccr1"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1347426666666666
$ws.Range("H7").Value = 0.404228
$ws.Range("I7").Value = 0.000968581456306409
$ws.Range("J7").Value = 0.0009685814563064089
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.979641
$ws.Range("N7").Value = 77.938923
$ws.Range("O7").Value = 0.3079672269977427
$ws.Range("P7").Value = 0.3079672269977428
$ws.Range("Q7").Value = 3.500566107382666
$ws.Range("R7").Value = 31.505094966444
$ws.Range("S7").Value = 0.0002982913452201201
$ws.Range("T7").Value = 0.0002982913452201201
